# REPORTGEN-709: update CWE full excel reports
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # Summary
$ws2 = $wb.Worksheets.Item(2)   # CWE (2011) Top25
$ws3 = $wb.Worksheets.Item(3)   # CWE (2011) Top25 details

# ---------------------------------------------------------------------------
# Sheet1 "Summary": turn the bare TECHNO_LOC / TECHNICAL_SIZING placeholders
# (rows 6 & 9) into proper two-column mini tables with header rows, using the
# already-blank rows 7 & 8 that sit between them.
# ---------------------------------------------------------------------------

# Row 6: header row for the Technology / Lines of Code mini table.
$ws1.Range("C6").Value = "Technology"
$ws1.Range("D6").Value = "Lines of Code"
$ws1.Range("B6").Copy()
$ws1.Range("C6:D6").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Row 7: the TECHNO_LOC table placeholder (used to live in C6) now lives here,
# with its header suppressed, plus a numeric-formatted value cell.
$ws1.Range("C7").Value = "RepGen:TABLE;TECHNO_LOC;HEADER=NO"
$ws1.Range("D7").NumberFormat = "0.00"

# Row 8: header row for the Characteristic / Value mini table.
$ws1.Range("C8").Value = "Characteristic"
$ws1.Range("D8").Value = "Value"
$ws1.Range("B6").Copy()
$ws1.Range("C8:D8").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Row 9: the TECHNICAL_SIZING table placeholder, header suppressed, plus a
# numeric-formatted value cell (N9/O9 stay put).
$ws1.Range("C9").Value = "RepGen:TABLE;TECHNICAL_SIZING;HEADER=NO"
$ws1.Range("D9").NumberFormat = "0.00"

# ---------------------------------------------------------------------------
# Sheet1 rows 12-14: findings summary table gains a real header row.
# ---------------------------------------------------------------------------
$ws1.Rows.Item(13).Insert()

$ws1.Range("B13").Value = "Quality Standard"
$ws1.Range("C13").Value = "Total Vulnerabilities"
$ws1.Range("D13").Value = "Added Vulnerabilities"
$ws1.Range("E13").Value = "Removed Vulnerabilities"
$ws1.Range("B12").Copy()
$ws1.Range("B13").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
$ws1.Range("B13").HorizontalAlignment = -4108  # xlCenter
$ws1.Range("B13").VerticalAlignment = -4108    # xlCenter
$ws1.Range("B13").WrapText = $true
$ws1.Range("C13:E13").Font.Name = "Verdana"
$ws1.Range("C13:E13").Font.Size = 9
$ws1.Range("C13:E13").Font.Bold = $true
$ws1.Range("C13:E13").Interior.Color = $ws1.Range("B12").Interior.Color
$ws1.Range("C13:E13").HorizontalAlignment = -4108  # xlCenter
$ws1.Range("C13:E13").VerticalAlignment = -4108    # xlCenter
$ws1.Range("C13:E13").WrapText = $true
$ws1.Rows.Item(13).RowHeight = 34.2

$ws1.Range("B14").Value = "RepGen:TABLE;QUALITY_STANDARDS_EVOLUTION;STD=CWE-2011-Top25,MORE=true,HEADER=NO"

# Column E is widened to fit the new "Removed Vulnerabilities" style header.
$ws1.Columns.Item(5).ColumnWidth = 18

# Nudge the logo shape so its cached bottom-right anchor is recomputed against
# the new column widths (purely cosmetic bookkeeping, position/size unchanged).
$shp = $ws1.Shapes.Item(1)
$shp.Left = $shp.Left

# ---------------------------------------------------------------------------
# Sheet2 "CWE (2011) Top25": becomes the rules/rationale table.
# ---------------------------------------------------------------------------
$ws2.Range("A1").Value = "Rules"
$ws2.Range("B1").Value = "Total Vulnerabilities"
$ws2.Range("C1").Value = "Added Vulnerabilities"
$ws2.Range("D1").Value = "Removed Vulnerabilities"
$ws2.Range("E1").Value = "Rationale"
$ws2.Range("F1").Value = "Description"
$ws2.Range("G1").Value = "Remediation"
$ws2.Range("A1").Copy()
$ws2.Range("B1:G1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws2.Range("A2").Value = "RepGen:TABLE;QUALITY_TAGS_RULES_EVOLUTION;STD=CWE-2011-Top25,DESC=true,HEADER=NO"
$ws2.Range("A2:G2").WrapText = $true
$ws2.Range("B2:D2").NumberFormat = "0.00"

$ws2.Columns.Item(2).ColumnWidth = 25.5546875
$ws2.Columns.Item(3).ColumnWidth = 24
$ws2.Columns.Item(4).ColumnWidth = 27.21875
$ws2.Columns.Item(5).ColumnWidth = 41.88671875
$ws2.Columns.Item(6).ColumnWidth = 36.5546875
$ws2.Columns.Item(7).ColumnWidth = 39.33203125

$ws2.Range("A2").Select()

# ---------------------------------------------------------------------------
# Sheet3 "CWE (2011) Top25 details": becomes the rule-violation detail table.
# ---------------------------------------------------------------------------
$ws3.Range("A1").Value = "Rule Name"
$ws3.Range("B1").Value = "Object Name"
$ws3.Range("C1").Value = "Object Type"
$ws3.Range("D1").Value = "Violation Status"
$ws3.Range("E1").Value = "Associated Value"
$ws3.Range("F1").Value = "File Path"
$ws3.Range("G1").Value = "Start Line"
$ws3.Range("H1").Value = "End Line"
$ws3.Range("A1").Copy()
$ws3.Range("B1:H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws3.Range("A2").Value = "RepGen:TABLE;LIST_RULES_VIOLATIONS_BOOKMARKS_TABLE;METRICS=CWE-2011-Top25,COUNT=-1,HEADER=NO"

$ws3.Columns.Item(1).ColumnWidth = 83
$ws3.Columns.Item(2).ColumnWidth = 58.5546875
$ws3.Columns.Item(3).ColumnWidth = 13.77734375
$ws3.Columns.Item(4).ColumnWidth = 18.21875
$ws3.Columns.Item(5).ColumnWidth = 19.21875
$ws3.Columns.Item(6).ColumnWidth = 60.5546875
$ws3.Columns.Item(7).ColumnWidth = 11.33203125
$ws3.Columns.Item(8).ColumnWidth = 10.109375

$ws3.Range("A3").Select()

$ws1.Activate()
